# Stage 1: update companies data
# Re-apply the new ordering/values for columns A,B,H,I,J,K across rows 3-11.
# Columns C,D,E,F,G (dates/status/source/time) stay attached to their row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = 'SEVEN (HOLDCO) LIMITED'
$c = $ws.Cells.Item(3, 2)
$c.NumberFormat = "@"
$c.Value = '16473606'
$c.Style = "Normal"
$ws.Cells.Item(3, 8).Value = 'Other'
$c = $ws.Cells.Item(3, 9)
$c.NumberFormat = "@"
$c.Value = '64209'
$c.Style = "Normal"
$ws.Cells.Item(3, 10).Value = 'Activities of other holding companies n.e.c.'
$ws.Cells.Item(3, 11).Value = 'Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles.'

$ws.Cells.Item(4, 1).Value = 'GANDER INVESTMENTS LTD'
$c = $ws.Cells.Item(4, 2)
$c.NumberFormat = "@"
$c.Value = '16473515'
$c.Style = "Normal"
$ws.Cells.Item(4, 8).Value = 'Investments'
$ws.Cells.Item(4, 9).Value = '68100,68209'
$ws.Cells.Item(4, 10).Value = ''
$ws.Cells.Item(4, 11).Value = ''

$ws.Cells.Item(5, 1).Value = 'INTERCONTINENTAL HOLDING COMPANY LIMITED'
$c = $ws.Cells.Item(5, 2)
$c.NumberFormat = "@"
$c.Value = '16473418'
$c.Style = "Normal"
$ws.Cells.Item(5, 8).Value = 'Other'
$c = $ws.Cells.Item(5, 9)
$c.NumberFormat = "@"
$c.Value = '64209'
$c.Style = "Normal"
$ws.Cells.Item(5, 10).Value = 'Activities of other holding companies n.e.c.'
$ws.Cells.Item(5, 11).Value = 'Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles.'

$ws.Cells.Item(6, 1).Value = 'TLJ INVESTMENT LTD'
$c = $ws.Cells.Item(6, 2)
$c.NumberFormat = "@"
$c.Value = '16473151'
$c.Style = "Normal"
$ws.Cells.Item(6, 8).Value = 'Investments'
$ws.Cells.Item(6, 9).Value = '41100,55100,68100'
$ws.Cells.Item(6, 10).Value = ''
$ws.Cells.Item(6, 11).Value = ''

$ws.Cells.Item(7, 1).Value = 'AJ INVESTMENT AND CONSULTANCY LTD'
$c = $ws.Cells.Item(7, 2)
$c.NumberFormat = "@"
$c.Value = '16473328'
$c.Style = "Normal"
$ws.Cells.Item(7, 8).Value = 'Investments'
$ws.Cells.Item(7, 9).Value = '64306,70229'
$ws.Cells.Item(7, 10).Value = 'Activities of real estate investment trusts'
$ws.Cells.Item(7, 11).Value = 'UK-regulated REIT companies.'

$ws.Cells.Item(8, 1).Value = 'GAUNT CAPITAL LTD'
$c = $ws.Cells.Item(8, 2)
$c.NumberFormat = "@"
$c.Value = '16473262'
$c.Style = "Normal"
$ws.Cells.Item(8, 8).Value = 'Capital'
$c = $ws.Cells.Item(8, 9)
$c.NumberFormat = "@"
$c.Value = '64209'
$c.Style = "Normal"
$ws.Cells.Item(8, 10).Value = 'Activities of other holding companies n.e.c.'
$ws.Cells.Item(8, 11).Value = 'Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles.'

$ws.Cells.Item(9, 1).Value = 'THE DISLEY GROUP LTD'
$c = $ws.Cells.Item(9, 2)
$c.NumberFormat = "@"
$c.Value = '16473398'
$c.Style = "Normal"
$ws.Cells.Item(9, 8).Value = 'Other'
$c = $ws.Cells.Item(9, 9)
$c.NumberFormat = "@"
$c.Value = '64209'
$c.Style = "Normal"
$ws.Cells.Item(9, 10).Value = 'Activities of other holding companies n.e.c.'
$ws.Cells.Item(9, 11).Value = 'Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles.'

$ws.Cells.Item(10, 1).Value = 'MARMIMI HOLDING LIMITED'
$c = $ws.Cells.Item(10, 2)
$c.NumberFormat = "@"
$c.Value = '16473234'
$c.Style = "Normal"
$ws.Cells.Item(10, 8).Value = 'Other'
$c = $ws.Cells.Item(10, 9)
$c.NumberFormat = "@"
$c.Value = '64209'
$c.Style = "Normal"
$ws.Cells.Item(10, 10).Value = 'Activities of other holding companies n.e.c.'
$ws.Cells.Item(10, 11).Value = 'Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles.'

$ws.Cells.Item(11, 1).Value = 'BRIDGEWICK PARTNERS LIMITED'
$c = $ws.Cells.Item(11, 2)
$c.NumberFormat = "@"
$c.Value = '16473142'
$c.Style = "Normal"
$ws.Cells.Item(11, 8).Value = 'Partners'
$c = $ws.Cells.Item(11, 9)
$c.NumberFormat = "@"
$c.Value = '64999'
$c.Style = "Normal"
$ws.Cells.Item(11, 10).Value = 'Financial intermediation not elsewhere classified'
$ws.Cells.Item(11, 11).Value = 'Catch-all credit-oriented SPVs for novel lending structures.'

